$d = $word.ActiveDocument

$d.Content.Find.Execute("This document contains the Project Definition and the Software Requirements Specification (SRS) for the mini-project: SmartTechMart — a retail application for a technology convenience store (Tech Store 4.0).", $false, $false, $false, $false, $false, $true, 1, $false, "This document contains the Project Definition and the Software Requirements Specification (SRS) for the mini-project: SmartTechMart — a retail application for a technology convenience store (Tech Store 4.0).", 2) | Out-Null
$d.Content.Find.Execute("SmartTechMart — Tech Convenience Store 4.0", $false, $false, $false, $false, $false, $true, 1, $false, "SmartTechMart — Tech Convenience Store 4.0", 2) | Out-Null
$d.Content.Find.Execute("Small technology retail and convenience stores increasingly need digital tools to manage stock, sales, customer loyalty, and online presence. SmartTechMart is a modern retail management application designed for small/medium tech convenience stores to handle point-of-sale (POS), inventory, supplier orders, basic e-commerce features, and analytics in a lightweight, affordable system.", $false, $false, $false, $false, $false, $true, 1, $false, "Small technology retail and convenience stores increasingly need digital tools to manage stock, sales, customer loyalty, and online presence. SmartTechMart is a modern retail management application designed for small/medium tech convenience stores to handle point-of-sale (POS), inventory, supplier orders, basic e-commerce features, and analytics in a lightweight, affordable system.", 2) | Out-Null
$d.Content.Find.Execute("- Provide a simple, reliable POS for in-store sales (barcode scanning, quick checkout).", $false, $false, $false, $false, $false, $true, 1, $false, "- Provide a simple, reliable POS for in-store sales (barcode scanning, quick checkout).", 2) | Out-Null
$d.Content.Find.Execute("- Enable online catalog browsing and click-and-collect for customers (optional MVP feature).", $false, $false, $false, $false, $false, $true, 1, $false, "- Enable online catalog browsing and click-and-collect for customers (optional MVP feature).", 2) | Out-Null
$d.Content.Find.Execute("- User accounts and role-based access (Admin, Manager, Cashier).", $false, $false, $false, $false, $false, $true, 1, $false, "- User accounts and role-based access (Admin, Manager, Cashier).", 2) | Out-Null
$d.Content.Find.Execute("- Advanced CRM features or marketing automation.", $false, $false, $false, $false, $false, $true, 1, $false, "- Advanced CRM features or marketing automation.", 2) | Out-Null
$d.Content.Find.Execute("- Store Manager: daily operations, inventory management, ordering.", $false, $false, $false, $false, $false, $true, 1, $false, "- Store Manager: daily operations, inventory management, ordering.", 2) | Out-Null
$d.Content.Find.Execute("- Budget and timeline constraints for a semester mini-project.", $false, $false, $false, $false, $false, $true, 1, $false, "- Budget and timeline constraints for a semester mini-project.", 2) | Out-Null
$d.Content.Find.Execute("- The store has an internet connection (some offline caching possible).", $false, $false, $false, $false, $false, $true, 1, $false, "- The store has an internet connection (some offline caching possible).", 2) | Out-Null
$d.Content.Find.Execute("- Basic working prototype: POS + Inventory + Reports.", $false, $false, $false, $false, $false, $true, 1, $false, "- Basic working prototype: POS + Inventory + Reports.", 2) | Out-Null
$d.Content.Find.Execute("- Week 2–3: Basic data model and inventory module implementation.", $false, $false, $false, $false, $false, $true, 1, $false, "- Week 2–3: Basic data model and inventory module implementation.", 2) | Out-Null
$d.Content.Find.Execute("This SRS describes the functional and non-functional requirements for SmartTechMart. It is intended for the development team, the project stakeholders, and the testers.", $false, $false, $false, $false, $false, $true, 1, $false, "This SRS describes the functional and non-functional requirements for SmartTechMart. It is intended for the development team, the project stakeholders, and the testers.", 2) | Out-Null
$d.Content.Find.Execute("Requirement IDs use the format FR-### for functional requirements and NFR-### for non-functional requirements.", $false, $false, $false, $false, $false, $true, 1, $false, "Requirement IDs use the format FR-### for functional requirements and NFR-### for non-functional requirements.", 2) | Out-Null
$d.Content.Find.Execute("2. Overall Description", $false, $false, $false, $false, $false, $true, 1, $false, "2. Overall Description", 2) | Out-Null
$d.Content.Find.Execute("SmartTechMart is a standalone retail management web application (progressive web app recommended) that can run on a local server or cloud host. It integrates with hardware devices such as barcode scanners and receipt printers, and optionally can integrate with payment terminals and simple QR payment gateways.", $false, $false, $false, $false, $false, $true, 1, $false, "SmartTechMart is a standalone retail management web application (progressive web app recommended) that can run on a local server or cloud host. It integrates with hardware devices such as barcode scanners and receipt printers, and optionally can integrate with payment terminals and simple QR payment gateways.", 2) | Out-Null
$d.Content.Find.Execute("- Support user accounts and role-based permissions.", $false, $false, $false, $false, $false, $true, 1, $false, "- Support user accounts and role-based permissions.", 2) | Out-Null
$d.Content.Find.Execute("- Manager: access to inventory, orders, reports.", $false, $false, $false, $false, $false, $true, 1, $false, "- Manager: access to inventory, orders, reports.", 2) | Out-Null
$d.Content.Find.Execute("Modern web browsers (Chrome, Edge, Firefox, Safari). Server: Linux or Windows hosting with Node.js / Python / PHP stack (team choice).", $false, $false, $false, $false, $false, $true, 1, $false, "Modern web browsers (Chrome, Edge, Firefox, Safari). Server: Linux or Windows hosting with Node.js / Python / PHP stack (team choice).", 2) | Out-Null
$d.Content.Find.Execute("- Availability of barcode scanner hardware and a network connection for cloud-hosted deployments.", $false, $false, $false, $false, $false, $true, 1, $false, "- Availability of barcode scanner hardware and a network connection for cloud-hosted deployments.", 2) | Out-Null
$d.Content.Find.Execute("Allow users to register (Admin only), login, reset password; support roles: Admin, Manager, Cashier. Priority: High.", $false, $false, $false, $false, $false, $true, 1, $false, "Allow users to register (Admin only), login, reset password; support roles: Admin, Manager, Cashier. Priority: High.", 2) | Out-Null
$d.Content.Find.Execute("Add, edit, delete products with fields: SKU, name, description, category, cost price, selling price, barcode, image, supplier. Priority: High.", $false, $false, $false, $false, $false, $true, 1, $false, "Add, edit, delete products with fields: SKU, name, description, category, cost price, selling price, barcode, image, supplier. Priority: High.", 2) | Out-Null
$d.Content.Find.Execute("Process sales with multiple items, apply discounts, accept payment type (cash/card/QR), generate receipt. Priority: High.", $false, $false, $false, $false, $false, $true, 1, $false, "Process sales with multiple items, apply discounts, accept payment type (cash/card/QR), generate receipt. Priority: High.", 2) | Out-Null
$d.Content.Find.Execute("Create supplier records and purchase orders, receive shipments to update inventory. Priority: Medium.", $false, $false, $false, $false, $false, $true, 1, $false, "Create supplier records and purchase orders, receive shipments to update inventory. Priority: Medium.", 2) | Out-Null
$d.Content.Find.Execute("Store customer info, track purchases, and apply loyalty points for discounts. Priority: Low/Optional.", $false, $false, $false, $false, $false, $true, 1, $false, "Store customer info, track purchases, and apply loyalty points for discounts. Priority: Low/Optional.", 2) | Out-Null
$d.Content.Find.Execute("Allow manual export and backup of database; import CSV for product bulk upload. Priority: Medium.", $false, $false, $false, $false, $false, $true, 1, $false, "Allow manual export and backup of database; import CSV for product bulk upload. Priority: Medium.", 2) | Out-Null
$d.Content.Find.Execute("Cache recent catalog and allow POS transactions offline with later sync. Priority: Low/Optional.", $false, $false, $false, $false, $false, $true, 1, $false, "Cache recent catalog and allow POS transactions offline with later sync. Priority: Low/Optional.", 2) | Out-Null
$d.Content.Find.Execute("System must handle up to 20 concurrent cashier sessions and respond to POS item scan within 1 second under typical load.", $false, $false, $false, $false, $false, $true, 1, $false, "System must handle up to 20 concurrent cashier sessions and respond to POS item scan within 1 second under typical load.", 2) | Out-Null
$d.Content.Find.Execute("Codebase should follow standard conventions, documented API, and automated tests for core flows.", $false, $false, $false, $false, $false, $true, 1, $false, "Codebase should follow standard conventions, documented API, and automated tests for core flows.", 2) | Out-Null
$d.Content.Find.Execute("Main entities: Product, Category, Supplier, Customer, User, Sale, SaleItem, PurchaseOrder, InventoryTransaction.", $false, $false, $false, $false, $false, $true, 1, $false, "Main entities: Product, Category, Supplier, Customer, User, Sale, SaleItem, PurchaseOrder, InventoryTransaction.", 2) | Out-Null
$d.Content.Find.Execute("UC-01: Process a Sale (Actor: Cashier)", $false, $false, $false, $false, $false, $true, 1, $false, "UC-01: Process a Sale (Actor: Cashier)", 2) | Out-Null
$d.Content.Find.Execute("4. Cashier selects payment type and records payment.", $false, $false, $false, $false, $false, $true, 1, $false, "4. Cashier selects payment type and records payment.", 2) | Out-Null
$d.Content.Find.Execute("Preconditions: Manager logged in.", $false, $false, $false, $false, $false, $true, 1, $false, "Preconditions: Manager logged in.", 2) | Out-Null
$d.Content.Find.Execute("Postconditions: Product available in catalog with initial stock (if provided).", $false, $false, $false, $false, $false, $true, 1, $false, "Postconditions: Product available in catalog with initial stock (if provided).", 2) | Out-Null
$d.Content.Find.Execute("2. On receipt, manager marks PO as received and inventory is updated.", $false, $false, $false, $false, $false, $true, 1, $false, "2. On receipt, manager marks PO as received and inventory is updated.", 2) | Out-Null
$d.Content.Find.Execute("- Inventory screens: filters, bulk import/export, low-stock highlights.", $false, $false, $false, $false, $false, $true, 1, $false, "- Inventory screens: filters, bulk import/export, low-stock highlights.", 2) | Out-Null
$d.Content.Find.Execute("- Core POS flow (sell, receipt, inventory update) works end-to-end in demo environment.", $false, $false, $false, $false, $false, $true, 1, $false, "- Core POS flow (sell, receipt, inventory update) works end-to-end in demo environment.", 2) | Out-Null
$d.Content.Find.Execute("8. Appendix", $false, $false, $false, $false, $false, $true, 1, $false, "8. Appendix", 2) | Out-Null
$lastPara = $d.Paragraphs.Last
$blank = $d.Paragraphs.Add($lastPara.Range)
$p1 = $d.Paragraphs.Add($blank.Range)
$p1.Range.InsertBefore("1")
$p2 = $d.Paragraphs.Add($p1.Range)
$p2.Range.InsertBefore("2")
$p3 = $d.Paragraphs.Add($p2.Range)
$p3.Range.InsertBefore("3")
$p4 = $d.Paragraphs.Add($p3.Range)
$p4.Range.InsertBefore("4")
$p5 = $d.Paragraphs.Add($p4.Range)
$p5.Range.InsertBefore("5")
